$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.836.64'
$ws.Range('E2').Value = '  -1.17%  '
$ws.Range('D3').Value = '1.662.63'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.33'
$ws.Range('E6').Value = '  +5.17%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +0.36%  '
$ws.Range('E9').Value = '  +0.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.19'
$ws.Range('E10').Value = '  +3.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0896'
$ws.Range('E11').Value = '  +3.82%  '
$ws.Range('D12').Value = '1.896.97'
$ws.Range('E12').Value = '  +0.16%  '
$ws.Range('D13').Value = '1.664.12'
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('E14').Value = '  +0.14%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.523'
$ws.Range('E15').Value = '  +0.55%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.95'
$ws.Range('E16').Value = '  +1.74%  '
$ws.Range('D17').Value = '26.836.30'
$ws.Range('E17').Value = '  -1.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '232.26'
$ws.Range('E18').Value = '  -2.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.89'
$ws.Range('E19').Value = '  +0.51%  '
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('E23').Value = '  -2.58%  '
$ws.Range('E24').Value = '  -1.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.80'
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.12'
$ws.Range('E26').Value = '  -0.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.116'
$ws.Range('E27').Value = '  +1.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.88'
$ws.Range('E28').Value = '  +0.39%  '
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('E30').Value = '  -0.32%  '
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('E32').Value = '  +1.74%  '
$ws.Range('D33').Value = '1.460.41'
$ws.Range('E33').Value = '  -5.13%  '
$ws.Range('E34').Value = '  +3.51%  '
$ws.Range('E35').Value = '  +2.98%  '
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.899'
$ws.Range('E37').Value = '  +1.32%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.573'
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('E39').Value = '  -0.28%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.79'
$ws.Range('E40').Value = '  -2.57%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('E42').Value = '  -0.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.977'
$ws.Range('E43').Value = '  +6.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '65.73'
$ws.Range('E44').Value = '  -0.76%  '
$ws.Range('D45').Value = '1.807.44'
$ws.Range('E45').Value = '  +0.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.776'
$ws.Range('E46').Value = '  +0.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.40'
$ws.Range('E47').Value = '  +0.53%  '
$ws.Range('E48').Value = '  +0.14%  '
$ws.Range('E49').Value = '  -1.35%  '
$ws.Range('E50').Value = '  +4.01%  '
$ws.Range('E51').Value = '  +0.36%  '
